$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns F, G, H
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy the header style (from existing header cell, e.g. A1) to new headers
$ws.Range("A1").Copy() | Out-Null
$ws.Range("F1:H1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Fill boolean data for rows 2-19, columns F, G, H - default FALSE
for ($r = 2; $r -le 19; $r++) {
    $ws.Cells.Item($r, 6).Value = $false
    $ws.Cells.Item($r, 7).Value = $false
    $ws.Cells.Item($r, 8).Value = $false
}

# Row 7 has RF_Outliers_MAD = TRUE
$ws.Cells.Item(7, 8).Value = $true
